# Add the new "python" / "programming language" user entry on the
# "Program" sheet (row 3), replacing the old "playwrightwithJava" /
# "Automationtool" values, and leave the selection on the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

$ws.Range("B3").Value = "python"
$ws.Range("C3").Value = "programming language"

$ws.Activate()
$ws.Range("C3").Select()
